$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218, shifting existing rows 218:266 down to 219:267
$ws.Rows.Item(218).Insert()

# Populate the newly-inserted row 218 with the new data record
$ws.Cells.Item(218, 1).Value = 3
$ws.Cells.Item(218, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(218, 3).Value = "Coquimbo"
$ws.Cells.Item(218, 4).Value = [DateTime]"2021-11-11"
$ws.Cells.Item(218, 5).Value = 5
$ws.Cells.Item(218, 6).Value = 100112017
$ws.Cells.Item(218, 7).Value = "Apio"
$ws.Cells.Item(218, 8).Value = "Americana (o)"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 160
$ws.Cells.Item(218, 11).Value = 9000
$ws.Cells.Item(218, 12).Value = 9000
$ws.Cells.Item(218, 13).Value = 9000
$ws.Cells.Item(218, 14).Value = "`$/docena de matas"
$ws.Cells.Item(218, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(218, 16).Value = 1500
$ws.Cells.Item(218, 17).Value = 6
$ws.Cells.Item(218, 18).Value = "Hortaliza"
